# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) previously held a "Strike#" style count; it is
# being regenerated to hold the actual strikeout (K) totals per game.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(1, 3, 9, 5, 3, 13, 7, 9, 11, 9, 6, 6, 8, 11, 5, 2, 7, 7, 8, 6, 8, 2, 6, 3, 3, 2)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
